$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Logs sheet: append the new mail-log entry as row 22
# ---------------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A22").Value = "Status van mijn bestelling"
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Wanneer wordt mijn bestelling bezorgd?"
$logs.Range("D22").Value = "Bestelling / Levering"
$logs.Range("F22").Value = "2025-06-22 18:52:13"
$logs.Range("G22").Value = "Nee"

# Extend the conditional-formatting ranges so the new row is covered too.
$logs.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))

# ---------------------------------------------------------------------------
# 2) Dashboard sheet: recomputed category counts (re-sorted) plus the new
#    "Bestelling / Levering" category in row 14
# ---------------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "IT / Technisch probleem"
$dash.Range("B2").Value = 4

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 3

$dash.Range("A4").Value = "Sollicitatie / Vacature"
$dash.Range("B4").Value = 2

$dash.Range("A5").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "Overig"
$dash.Range("B6").Value = 2

$dash.Range("A7").Value = "Openingstijden / Locatie"
$dash.Range("B7").Value = 1

$dash.Range("A8").Value = "Uitnodiging / Evenement"
$dash.Range("B8").Value = 1

$dash.Range("A9").Value = "Klacht / Probleem"
$dash.Range("B9").Value = 1

$dash.Range("A10").Value = "Offerte / Prijsaanvraag"
$dash.Range("B10").Value = 1

$dash.Range("A11").Value = "Retour / Terugbetaling"
$dash.Range("B11").Value = 1

$dash.Range("A12").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B12").Value = 1

$dash.Range("A13").Value = "Factuur / Administratie"
$dash.Range("B13").Value = 1

$dash.Range("A14").Value = "Bestelling / Levering"
$dash.Range("B14").Value = 1

# ---------------------------------------------------------------------------
# 3) Chart on the Dashboard sheet: extend the series references to row 14
# ---------------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$14,'Dashboard'!`$B`$2:`$B`$14,1)"
